$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The JS-implementation sweep reached a few more rows since the last save:
# NotificationApplicationServiceTest.java (rows 24-26), GroupResourceTest.java
# (row 95) and the *ResourceTest.java block (rows 100-106) are now done too,
# so flag column D ("JS Implemented") with "Y" for each of them. The D107/E107
# summary formulas (COUNTIF / ratio) recalc automatically from this.
$doneRows = @(24, 25, 26, 95, 100, 101, 102, 103, 104, 105, 106)
foreach ($r in $doneRows) {
    $ws.Range("D$r").Value = "Y"
}

# Leave the sheet scrolled to where the work left off and the cursor on the
# next unmarked row (D96) instead of the top of the list.
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 1
$ws.Range("D96").Select()
